$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column B keeps its text format so long numeric-looking
# terminal IDs are not coerced into floating point numbers.
$ws.Range("B2:B46").NumberFormat = "@"

$ws.Range("A2").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B2").Value = '202号直流'
$ws.Range("C2").Value = 46046.067511574074
$ws.Range("D2").Value = 46048.22355324074

$ws.Range("A3").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B3").Value = '212号直流'
$ws.Range("C3").Value = 46046.751238425924
$ws.Range("D3").Value = 46048.22355324074

$ws.Range("A4").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B4").Value = '204号直流'
$ws.Range("C4").Value = 46047.409224537034
$ws.Range("D4").Value = 46048.22355324074

$ws.Range("A5").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B5").Value = '302号直流'
$ws.Range("C5").Value = 46047.427141203705
$ws.Range("D5").Value = 46048.22355324074

$ws.Range("A6").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B6").Value = '203号直流'
$ws.Range("C6").Value = 46047.527025462965
$ws.Range("D6").Value = 46048.22355324074

$ws.Range("A7").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B7").Value = '301号直流'
$ws.Range("C7").Value = 46047.55606481482
$ws.Range("D7").Value = 46048.22355324074

$ws.Range("A8").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B8").Value = '112号直流'
$ws.Range("C8").Value = 46047.56303240741
$ws.Range("D8").Value = 46048.22355324074

$ws.Range("A9").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B9").Value = '205号直流'
$ws.Range("C9").Value = 46047.57376157407
$ws.Range("D9").Value = 46048.22355324074

$ws.Range("A10").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B10").Value = '105号直流'
$ws.Range("C10").Value = 46047.58038194444
$ws.Range("D10").Value = 46048.22355324074

$ws.Range("A11").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B11").Value = '108号直流'
$ws.Range("C11").Value = 46047.582766203705
$ws.Range("D11").Value = 46048.22355324074

$ws.Range("A12").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B12").Value = '104号直流'
$ws.Range("C12").Value = 46047.67915509259
$ws.Range("D12").Value = 46048.22355324074

$ws.Range("A13").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B13").Value = '208号直流'
$ws.Range("C13").Value = 46047.69856481482
$ws.Range("D13").Value = 46048.22355324074

$ws.Range("A14").Value = '长沙市开福区高岭香江国际城充电站建设项目'
$ws.Range("B14").Value = '102号直流'
$ws.Range("C14").Value = 46047.70590277778
$ws.Range("D14").Value = 46048.22355324074

$ws.Range("A15").Value = '飞狐四方坪西区充电站'
$ws.Range("B15").Value = '9176699400500102'
$ws.Range("C15").Value = 46044.55170138889
$ws.Range("D15").Value = 46048.312569444446

$ws.Range("A16").Value = '飞狐四方坪西区充电站'
$ws.Range("B16").Value = '9176699400500804'
$ws.Range("C16").Value = 46046.038125
$ws.Range("D16").Value = 46048.312569444446

$ws.Range("A17").Value = '飞狐四方坪西区充电站'
$ws.Range("B17").Value = '9176699400501303'
$ws.Range("C17").Value = 46046.5399537037
$ws.Range("D17").Value = 46048.312569444446

$ws.Range("A18").Value = '飞狐四方坪西区充电站'
$ws.Range("B18").Value = '9176699400500304'
$ws.Range("C18").Value = 46046.547060185185
$ws.Range("D18").Value = 46048.312569444446

$ws.Range("A19").Value = '飞狐四方坪西区充电站'
$ws.Range("B19").Value = '9176699400501203'
$ws.Range("C19").Value = 46046.54938657407
$ws.Range("D19").Value = 46048.312569444446

$ws.Range("A20").Value = '飞狐四方坪西区充电站'
$ws.Range("B20").Value = '9176699400500205'
$ws.Range("C20").Value = 46046.57809027778
$ws.Range("D20").Value = 46048.312569444446

$ws.Range("A21").Value = '飞狐四方坪西区充电站'
$ws.Range("B21").Value = '9176699400500302'
$ws.Range("C21").Value = 46046.58342592593
$ws.Range("D21").Value = 46048.312569444446

$ws.Range("A22").Value = '飞狐四方坪西区充电站'
$ws.Range("B22").Value = '9176699400500403'
$ws.Range("C22").Value = 46046.70376157408
$ws.Range("D22").Value = 46048.312569444446

$ws.Range("A23").Value = '飞狐四方坪东区充电站'
$ws.Range("B23").Value = '9176699442100402'
$ws.Range("C23").Value = 46046.92711805556
$ws.Range("D23").Value = 46048.312569444446

$ws.Range("A24").Value = '飞狐四方坪南区充电站'
$ws.Range("B24").Value = '9176699368200101'
$ws.Range("C24").Value = 46047.03225694445
$ws.Range("D24").Value = 46048.312569444446

$ws.Range("A25").Value = '飞狐四方坪东区充电站'
$ws.Range("B25").Value = '9176699442100101'
$ws.Range("C25").Value = 46047.11447916667
$ws.Range("D25").Value = 46048.312569444446

$ws.Range("A26").Value = '飞狐四方坪西区充电站'
$ws.Range("B26").Value = '9176699400501102'
$ws.Range("C26").Value = 46047.412777777776
$ws.Range("D26").Value = 46048.312569444446

$ws.Range("A27").Value = '飞狐四方坪西区充电站'
$ws.Range("B27").Value = '9176699400500203'
$ws.Range("C27").Value = 46047.49927083333
$ws.Range("D27").Value = 46048.312569444446

$ws.Range("A28").Value = '飞狐四方坪东区充电站'
$ws.Range("B28").Value = '9176699425700301'
$ws.Range("C28").Value = 46047.515694444446
$ws.Range("D28").Value = 46048.312569444446

$ws.Range("A29").Value = '飞狐四方坪西区充电站'
$ws.Range("B29").Value = '9176699400500405'
$ws.Range("C29").Value = 46047.51994212963
$ws.Range("D29").Value = 46048.312569444446

$ws.Range("A30").Value = '飞狐四方坪南区充电站'
$ws.Range("B30").Value = '9176699368200306'
$ws.Range("C30").Value = 46047.52872685185
$ws.Range("D30").Value = 46048.312569444446

$ws.Range("A31").Value = '飞狐四方坪西区充电站'
$ws.Range("B31").Value = '9176699400501205'
$ws.Range("C31").Value = 46047.534791666665
$ws.Range("D31").Value = 46048.312569444446

$ws.Range("A32").Value = '飞狐四方坪南区充电站'
$ws.Range("B32").Value = '9176699368200202'
$ws.Range("C32").Value = 46047.54965277778
$ws.Range("D32").Value = 46048.312569444446

$ws.Range("A33").Value = '飞狐四方坪西区充电站'
$ws.Range("B33").Value = '9176699400500601'
$ws.Range("C33").Value = 46047.55212962963
$ws.Range("D33").Value = 46048.312569444446

$ws.Range("A34").Value = '飞狐四方坪西区充电站'
$ws.Range("B34").Value = '9176699400500502'
$ws.Range("C34").Value = 46047.55708333333
$ws.Range("D34").Value = 46048.312569444446

$ws.Range("A35").Value = '飞狐四方坪东区充电站'
$ws.Range("B35").Value = '9176699435600102'
$ws.Range("C35").Value = 46047.563576388886
$ws.Range("D35").Value = 46048.312569444446

$ws.Range("A36").Value = '飞狐四方坪西区充电站'
$ws.Range("B36").Value = '9176699400500501'
$ws.Range("C36").Value = 46047.578043981484
$ws.Range("D36").Value = 46048.312569444446

$ws.Range("A37").Value = '飞狐四方坪南区充电站'
$ws.Range("B37").Value = '9176699368200103'
$ws.Range("C37").Value = 46047.58162037037
$ws.Range("D37").Value = 46048.312569444446

$ws.Range("A38").Value = '飞狐四方坪西区充电站'
$ws.Range("B38").Value = '9176699400500104'
$ws.Range("C38").Value = 46047.59732638889
$ws.Range("D38").Value = 46048.312569444446

$ws.Range("A39").Value = '飞狐四方坪东区充电站'
$ws.Range("B39").Value = '9176699442100201'
$ws.Range("C39").Value = 46047.59778935185
$ws.Range("D39").Value = 46048.312569444446

$ws.Range("A40").Value = '飞狐四方坪南区充电站'
$ws.Range("B40").Value = '9176699368200406'
$ws.Range("C40").Value = 46047.60381944444
$ws.Range("D40").Value = 46048.312569444446

$ws.Range("A41").Value = '飞狐四方坪西区充电站'
$ws.Range("B41").Value = '9176699400500204'
$ws.Range("C41").Value = 46047.61420138889
$ws.Range("D41").Value = 46048.312569444446

$ws.Range("A42").Value = '飞狐四方坪西区充电站'
$ws.Range("B42").Value = '9176699355900102'
$ws.Range("C42").Value = 46047.64368055556
$ws.Range("D42").Value = 46048.312569444446

$ws.Range("A43").Value = '飞狐四方坪东区充电站'
$ws.Range("B43").Value = '9176699442100802'
$ws.Range("C43").Value = 46047.6596412037
$ws.Range("D43").Value = 46048.312569444446

$ws.Range("A44").Value = '飞狐四方坪西区充电站'
$ws.Range("B44").Value = '9176699400500802'
$ws.Range("C44").Value = 46047.6925
$ws.Range("D44").Value = 46048.312569444446

$ws.Range("A45").Value = '飞狐四方坪西区充电站'
$ws.Range("B45").Value = '9176699400501101'
$ws.Range("C45").Value = 46047.70997685185
$ws.Range("D45").Value = 46048.312569444446

$ws.Range("A46").Value = '飞狐四方坪西区充电站'
$ws.Range("B46").Value = '9176699400500202'
$ws.Range("C46").Value = 46047.7109837963
$ws.Range("D46").Value = 46048.312569444446

$ws.Range("E21").Select()